# Add two new PSP log entries (2019-12-05 .. 2019-12-07), replacing the
# previous placeholder/empty rows 39-42 with five populated rows (39-43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Seed formatting for the five data rows by copying from existing,
#    already-correctly-formatted rows, then overwrite with the new values.
#    - Column A needs the date number format used by row 6 (m"월" d"일").
#    - Columns B:F need the border/number-format treatment already used
#      by row 41 (time columns, integer columns, text column).
# ---------------------------------------------------------------------
$ws.Range("A6").Copy($ws.Range("A39"))
$ws.Range("A6").Copy($ws.Range("A40"))
$ws.Range("A6").Copy($ws.Range("A41"))
$ws.Range("A6").Copy($ws.Range("A42"))
$ws.Range("A6").Copy($ws.Range("A43"))

$ws.Range("B41:F41").Copy($ws.Range("B39:F39"))
$ws.Range("B41:F41").Copy($ws.Range("B40:F40"))
$ws.Range("B41:F41").Copy($ws.Range("B42:F42"))
$ws.Range("B41:F41").Copy($ws.Range("B43:F43"))
# B41:F41 already holds its own formatting.

# Row heights to roughly match the source rows.
$ws.Rows(39).RowHeight = 15
$ws.Rows(40).RowHeight = 15
$ws.Rows(43).RowHeight = 13.5

# Columns D/E in the new rows display as plain (General) centered numbers
# rather than the bordered "0" integer format used elsewhere.
$ws.Range("D39:E43").NumberFormat = "General"
$ws.Range("D39:E43").HorizontalAlignment = -4108  # xlCenter

# ---------------------------------------------------------------------
# 2) Write the actual log values.
# ---------------------------------------------------------------------

# Row 39: Thu 2019-12-05, 21:30 -> 01:00 (+1d), 0 interrupt, 210 delta
$ws.Range("A39").Value = 43804
$ws.Range("B39").Value = 21.5/24
$ws.Range("C39").Value = 1/24
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 210

# Row 40: Fri 2019-12-06, 11:00 -> 18:00, 30 interrupt, 270 delta
$ws.Range("A40").Value = 43805
$ws.Range("B40").Value = 11/24
$ws.Range("C40").Value = 18/24
$ws.Range("D40").Value = 30
$ws.Range("E40").Value = 270

# Row 41: Fri 2019-12-06, 23:00 -> 03:30 (+1d), 30 interrupt, 180 delta
$ws.Range("A41").Value = 43805
$ws.Range("B41").Value = 23/24
$ws.Range("C41").Value = 3.5/24
$ws.Range("D41").Value = 30
$ws.Range("E41").Value = 180

# Row 42: Sat 2019-12-07, 11:00 -> 13:00, 0 interrupt, 120 delta
$ws.Range("A42").Value = 43806
$ws.Range("B42").Value = 11/24
$ws.Range("C42").Value = 13/24
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 120

# Row 43: Sat 2019-12-07, 15:00 -> 23:00, 90 interrupt, 390 delta
$ws.Range("A43").Value = 43806
$ws.Range("B43").Value = 15/24
$ws.Range("C43").Value = 23/24
$ws.Range("D43").Value = 90
$ws.Range("E43").Value = 390

# ---------------------------------------------------------------------
# 3) Activity text (column F) — two distinct log descriptions, each
#    written once as rich text with "개발" in a distinct run (matching
#    the author's formatting), then copied into the repeated cells so
#    the three/two occurrences share one underlying string entry.
# ---------------------------------------------------------------------
$devPersonal = "개발 - 개인 요람 생성 및 맞춤과목 리스트 구현"
$devTimetable = "개발 - 시간표 데이터 수치화"

$ws.Range("F39").Value = $devPersonal
$ws.Range("F39").Characters(1,2).Font.Name = "맑은 고딕"
$ws.Range("F39").Copy($ws.Range("F40"))
$ws.Range("F39").Copy($ws.Range("F43"))

$ws.Range("F41").Value = $devTimetable
$ws.Range("F41").Characters(3, $devTimetable.Length-2).Font.Name = "돋움"
$ws.Range("F41").Copy($ws.Range("F42"))

# ---------------------------------------------------------------------
# 4) Sheet-level bookkeeping to mirror the saved view state.
# ---------------------------------------------------------------------
$ws.Range("F47").Select()
